$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

Set-TextValue 'D2' '28.039.96'
Set-TextValue 'E2' '  +2.23%  '
Set-TextValue 'D3' '1.910.55'
Set-TextValue 'E3' '  +2.50%  '
Set-TextValue 'D4' '1.005'
Set-TextValue 'E4' '  -0.74%  '
Set-TextValue 'D5' '315.38'
Set-TextValue 'E5' '  +1.32%  '
Set-TextValue 'E6' '  -0.80%  '
Set-TextValue 'D7' '0.4804'
Set-TextValue 'E7' '  +0.71%  '
Set-TextValue 'D8' '0.3809'
Set-TextValue 'E8' '  +1.43%  '
Set-TextValue 'D9' '0.07358'
Set-TextValue 'E9' '  +0.57%  '
Set-TextValue 'D10' '0.9332'
Set-TextValue 'E10' '  -0.06%  '
Set-TextValue 'D11' '20.83'
Set-TextValue 'E11' '  +0.92%  '
Set-TextValue 'D12' '0.07776'
Set-TextValue 'D13' '1.882.43'
Set-TextValue 'E13' '  +0.61%  '
Set-TextValue 'D14' '5.493'
Set-TextValue 'D15' '6.634'
Set-TextValue 'E15' '  +1.25%  '
Set-TextValue 'D16' '92.10'
Set-TextValue 'E16' '  +1.73%  '
Set-TextValue 'E17' '  -0.81%  '
Set-TextValue 'D18' '0.000008862'
Set-TextValue 'E18' '  -0.12%  '
Set-TextValue 'E19' '  -0.73%  '
Set-TextValue 'D20' '28.069.96'
Set-TextValue 'E20' '  +2.20%  '
Set-TextValue 'E21' '  +0.75%  '
Set-TextValue 'D22' '5.163'
Set-TextValue 'E22' '  +1.04%  '
Set-TextValue 'D23' '2.125.62'
Set-TextValue 'E23' '  +0.76%  '
Set-TextValue 'E24' '  +2.16%  '
Set-TextValue 'D25' '155.56'
Set-TextValue 'E25' '  +0.18%  '
Set-TextValue 'E26' '  -1.07%  '
Set-TextValue 'D27' '18.48'
Set-TextValue 'E27' '  +0.10%  '
Set-TextValue 'D28' '2.137'
Set-TextValue 'E28' '  +6.01%  '
Set-TextValue 'D29' '116.91'
Set-TextValue 'E29' '  +1.43%  '
Set-TextValue 'D30' '4.963'
Set-TextValue 'E30' '  +0.58%  '
Set-TextValue 'D31' '0.08944'
Set-TextValue 'E31' '  +0.60%  '
Set-TextValue 'D32' '3.304'
Set-TextValue 'E32' '  -0.59%  '
Set-TextValue 'D33' '1.263'
Set-TextValue 'E33' '  +4.10%  '
Set-TextValue 'D34' '0.7773'
Set-TextValue 'E34' '  +3.26%  '
Set-TextValue 'D35' '4.675'
Set-TextValue 'E35' '  +1.74%  '
Set-TextValue 'D36' '2.629'
Set-TextValue 'E36' '  -3.64%  '
Set-TextValue 'D37' '0.02054'
Set-TextValue 'E37' '  +0.80%  '
Set-TextValue 'D38' '1.112'
Set-TextValue 'E38' '  -0.33%  '
Set-TextValue 'D39' '0.05312'
Set-TextValue 'B40' 'MXToken'
Set-TextValue 'C40' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue 'D40' '3.005'
Set-TextValue 'E40' '  +0.67%  '
Set-TextValue 'B41' 'TheSandbox'
Set-TextValue 'C41' 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
Set-TextValue 'D41' '0.5488'
Set-TextValue 'E41' '  +3.18%  '
Set-TextValue 'E42' '  -0.65%  '
Set-TextValue 'D43' '0.1529'
Set-TextValue 'E43' '  +0.48%  '
Set-TextValue 'D44' '8.481'
Set-TextValue 'E44' '  -1.00%  '
Set-TextValue 'D45' '10.72'
Set-TextValue 'E45' '  +1.19%  '
Set-TextValue 'D46' '0.4825'
Set-TextValue 'E46' '  +0.61%  '
Set-TextValue 'D47' '108.18'
Set-TextValue 'E47' '  +5.24%  '
Set-TextValue 'E48' '  -0.87%  '
Set-TextValue 'D49' '1.649'
Set-TextValue 'E49' '  -0.35%  '
Set-TextValue 'D50' '67.87'
Set-TextValue 'E50' '  +1.03%  '
Set-TextValue 'D51' '0.06071'
Set-TextValue 'E51' '  -0.11%  '
